# Apply "Apio" (Apio - Mercado Mayorista Lo Valledor de Santiago) weekly update:
# - Existing weekly records (rows 662-711) shift down by one pair (Primera/Segunda),
#   i.e. each pair of rows takes on the Fecha/Volumen/Precios/Origen of the pair that
#   preceded it.
# - A brand-new pair of records is inserted at the top (rows 662-663).
# - The last existing pair is duplicated onto two new rows appended at the bottom
#   (712-713), growing the sheet from A1:R711 to A1:R713.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target row, Calidad (Primera/Segunda), Fecha (serial), Volumen, Precio
# minimo/maximo/promedio, Origen, Precio $/Kg. All other columns (Mercado ID, Mercado,
# Region, Codreg, Categoria ID/Categoria/Variedad, Unidad de comercializacion, Kg o
# Unidades, Clasificacion) are identical for every row in this block and are left as-is
# for the rows that already existed (662-711); they are written explicitly for the two
# brand-new rows (712-713).
$data = @(
  @{ Row=662; I="Primera"; D=44746; J=2600; K=7000; L=8000; M=7538; O="Región de Coquimbo"; P=1256 },
  @{ Row=663; I="Segunda"; D=44746; J=600; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=664; I="Primera"; D=44386; J=1400; K=6000; L=7000; M=6429; O="Región de Coquimbo"; P=1072 },
  @{ Row=665; I="Segunda"; D=44386; J=760; K=3000; L=4000; M=3658; O="Región de Coquimbo"; P=610 },
  @{ Row=666; I="Primera"; D=44690; J=2100; K=7000; L=8000; M=7429; O="Región de Coquimbo"; P=1238 },
  @{ Row=667; I="Segunda"; D=44690; J=900; K=5000; L=6000; M=5667; O="Región de Coquimbo"; P=944 },
  @{ Row=668; I="Primera"; D=44631; J=1800; K=8000; L=9000; M=8333; O="Región de Coquimbo"; P=1389 },
  @{ Row=669; I="Segunda"; D=44631; J=900; K=5000; L=6000; M=5667; O="Región de Coquimbo"; P=944 },
  @{ Row=670; I="Primera"; D=44235; J=970; K=5000; L=6000; M=5500; O="Región de Coquimbo"; P=917 },
  @{ Row=671; I="Segunda"; D=44235; J=340; K=4000; L=4000; M=4000; O="Región de Coquimbo"; P=667 },
  @{ Row=672; I="Primera"; D=44582; J=2100; K=5500; L=6000; M=5714; O="Región de Coquimbo"; P=952 },
  @{ Row=673; I="Segunda"; D=44582; J=900; K=4000; L=4500; M=4333; O="Región de Coquimbo"; P=722 },
  @{ Row=674; I="Primera"; D=44307; J=1400; K=5000; L=6000; M=5536; O="Región de Coquimbo"; P=923 },
  @{ Row=675; I="Segunda"; D=44307; J=600; K=3000; L=4000; M=3583; O="Región de Coquimbo"; P=597 },
  @{ Row=676; I="Primera"; D=44672; J=2600; K=7000; L=8000; M=7538; O="Región de Coquimbo"; P=1256 },
  @{ Row=677; I="Segunda"; D=44672; J=1800; K=5000; L=6000; M=5667; O="Región de Coquimbo"; P=944 },
  @{ Row=678; I="Primera"; D=44344; J=1500; K=7000; L=8000; M=7400; O="Región de Coquimbo"; P=1233 },
  @{ Row=679; I="Segunda"; D=44344; J=600; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=680; I="Primera"; D=44433; J=2000; K=7000; L=8000; M=7600; O="Región de Coquimbo"; P=1267 },
  @{ Row=681; I="Segunda"; D=44433; J=700; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=682; I="Primera"; D=44707; J=1900; K=6000; L=7000; M=6368; O="Región de Coquimbo"; P=1061 },
  @{ Row=683; I="Segunda"; D=44707; J=800; K=4000; L=5000; M=4625; O="Región de Coquimbo"; P=771 },
  @{ Row=684; I="Primera"; D=44265; J=520; K=6000; L=7000; M=6519; O="Provincia de Limarí"; P=1086 },
  @{ Row=685; I="Segunda"; D=44265; J=180; K=4000; L=4000; M=4000; O="Provincia de Limarí"; P=667 },
  @{ Row=686; I="Primera"; D=44421; J=1100; K=7000; L=8000; M=7545; O="Región de Coquimbo"; P=1258 },
  @{ Row=687; I="Segunda"; D=44421; J=400; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=688; I="Primera"; D=44215; J=520; K=7000; L=8000; M=7500; O="Región de Coquimbo"; P=1250 },
  @{ Row=689; I="Segunda"; D=44215; J=210; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=690; I="Primera"; D=44566; J=2000; K=7000; L=8000; M=7600; O="Región de Coquimbo"; P=1267 },
  @{ Row=691; I="Segunda"; D=44566; J=600; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=692; I="Primera"; D=44637; J=2600; K=8000; L=9000; M=8462; O="Región de Coquimbo"; P=1410 },
  @{ Row=693; I="Segunda"; D=44637; J=1100; K=5000; L=6000; M=5727; O="Región de Coquimbo"; P=954 },
  @{ Row=694; I="Primera"; D=44195; J=2100; K=7000; L=8000; M=7571; O="Región de Coquimbo"; P=1262 },
  @{ Row=695; I="Segunda"; D=44195; J=800; K=5000; L=5000; M=5000; O="Región de Coquimbo"; P=833 },
  @{ Row=696; I="Primera"; D=44244; J=1500; K=7000; L=8000; M=7600; O="Región de Coquimbo"; P=1267 },
  @{ Row=697; I="Segunda"; D=44244; J=600; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=698; I="Primera"; D=44442; J=1400; K=7000; L=8000; M=7486; O="Región de Coquimbo"; P=1248 },
  @{ Row=699; I="Segunda"; D=44442; J=570; K=6000; L=6000; M=6000; O="Región de Coquimbo"; P=1000 },
  @{ Row=700; I="Primera"; D=44483; J=5400; K=6000; L=7000; M=6519; O="Región de Coquimbo"; P=1086 },
  @{ Row=701; I="Segunda"; D=44483; J=1400; K=4000; L=5000; M=4571; O="Región de Coquimbo"; P=762 },
  @{ Row=702; I="Primera"; D=44663; J=1430; K=9000; L=10000; M=9476; O="Provincia del Elquí"; P=1579 },
  @{ Row=703; I="Segunda"; D=44663; J=470; K=8000; L=8000; M=8000; O="Provincia del Elquí"; P=1333 },
  @{ Row=704; I="Primera"; D=44188; J=2000; K=5500; L=6000; M=5762; O="Región de Coquimbo"; P=960 },
  @{ Row=705; I="Segunda"; D=44188; J=600; K=4000; L=4000; M=4000; O="Región de Coquimbo"; P=667 },
  @{ Row=706; I="Primera"; D=44187; J=2000; K=5500; L=6000; M=5700; O="Región de Coquimbo"; P=950 },
  @{ Row=707; I="Segunda"; D=44187; J=600; K=4500; L=4500; M=4500; O="Región de Coquimbo"; P=750 },
  @{ Row=708; I="Primera"; D=44519; J=1370; K=5000; L=6000; M=5562; O="Región de Coquimbo"; P=927 },
  @{ Row=709; I="Segunda"; D=44519; J=560; K=4000; L=4000; M=4000; O="Región de Coquimbo"; P=667 },
  @{ Row=710; I="Primera"; D=44231; J=880; K=5000; L=6000; M=5500; O="Región de Coquimbo"; P=917 },
  @{ Row=711; I="Segunda"; D=44231; J=250; K=4000; L=4000; M=4000; O="Región de Coquimbo"; P=667 },
  @{ Row=712; I="Primera"; D=44194; J=1500; K=6000; L=7000; M=6533; O="Región de Coquimbo"; P=1089 },
  @{ Row=713; I="Segunda"; D=44194; J=500; K=5000; L=5000; M=5000; O="Región de Coquimbo"; P=833 }
)

foreach ($r in $data) {
  $ws.Range("D" + $r.Row).Value = $r.D
  $ws.Range("J" + $r.Row).Value = $r.J
  $ws.Range("K" + $r.Row).Value = $r.K
  $ws.Range("L" + $r.Row).Value = $r.L
  $ws.Range("M" + $r.Row).Value = $r.M
  $ws.Range("O" + $r.Row).Value = $r.O
  $ws.Range("P" + $r.Row).Value = $r.P
}

# The two brand-new rows (712-713) did not exist before, so fill in the remaining
# (previously-constant) columns explicitly, copying them from the row above.
$lastCols = @("A","B","C","E","F","G","H","I","N","Q","R")
foreach ($col in $lastCols) {
  $ws.Range($col + "712").Value = $ws.Range($col + "710").Value2
  $ws.Range($col + "713").Value = $ws.Range($col + "711").Value2
}

# Match the date-formatted style used by the rest of column D.
$ws.Range("D712").NumberFormat = $ws.Range("D710").NumberFormat
$ws.Range("D713").NumberFormat = $ws.Range("D711").NumberFormat
